# How to Write an Academic Paper.docx -- apply the "complete" commit edits.
# Each call does a targeted Find & Replace so only the edited words/phrases
# change, leaving the rest of each paragraph's runs untouched.

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

# --- 2.2 topic-selection bullets ---------------------------------------
Replace-Text "deciding a topic" "choosing a topic"
Replace-Text "written on some many times" "written for so many times"
Replace-Text "we have a unique approach to the topic" "we have a unique and impressive approach to the topic"

Replace-Text "Read articles, journals and papers" "Read articles, books and papers"
Replace-Text "finding words commonly used" "finding keywords commonly used"

Replace-Text "reading related works/topics in encyclopedias, articles" "reading related works and topics in papers, articles"

Replace-Text "make sure that choose an area we can afford" "make sure that choose one that we can afford"

# --- 3.2 paper proposal bullets -----------------------------------------
Replace-Text "knowledge of the area and helps" "knowledge of the problem and helps"
Replace-Text "other scientific areas." "other scientific fields."

Replace-Text "stated in the introductory paragraph." "stated in the introduction chapter."

# --- 3.3 outline section --------------------------------------------------
Replace-Text "any paper. Outlining before writing the rough draft" "any paper. Writing an outline before writing the rough draft"
Replace-Text "information is considered supporting ideas. Besides" "information is considered supporting our ideas. Besides"

# --- Introduction bullet --------------------------------------------------
Replace-Text "summarizes background data on the topic" "summarizes background information on the topic"
Replace-Text "providing the main goal of the work." "providing the main goal of our work."

# --- Our work and contribution bullet --------------------------------------
Replace-Text "as well as how we decide to go about finding the answers, and we should" "as well as how we find the solution, and we should"

# --- Experiment or evaluation bullet ---------------------------------------
Replace-Text "thoroughly describe the results we obtained." "thoroughly describe the experiment results we obtained."
Replace-Text "pictorially using graphs or histograms. In addition" "pictorially using figures and tables. In addition"
Replace-Text "tell the reader what all these data mean." "tell the readers what these data mean."

# --- References bullet ------------------------------------------------------
Replace-Text "is a Reference section listing" "is a reference section listing"

# --- 3.4 Fill in Each Chapter ------------------------------------------------
Replace-Text "we can try to express in these ways." "we can try to write it in these ways."

# --- 3.5 Revise the Paper Iteratively bullets --------------------------------
Replace-Text "support for our argument" "support for our ideas"

# --- Summary section ---------------------------------------------------------
Replace-Text "research project in CS or SE." "research project in Computer Science or Software Engineering."
Replace-Text "Write the program according to the problem." "Write the program to solve the problem."

# --- Drop the stale lastRenderedPageBreak cache hint before the Jack
#     Caulfield reference: round-tripping the run's text through Find &
#     Replace re-creates the run without the non-text child element.
Replace-Text "Jack Caulfield. How to Write an Essay Outline. " "Jack Caulfield. How to Write an Essay Outline. "

Write-Output "done"
